$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.828.24'
$ws.Range('E2').Value = '  -0.70%  '
$ws.Range('D3').Value = '3.366.12'
$ws.Range('E3').Value = '  +1.41%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '403.18'
$ws.Range('E5').Value = '  -1.55%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '127.57'
$ws.Range('E6').Value = '  +12.78%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.604'
$ws.Range('E7').Value = '  +6.65%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.670'
$ws.Range('E9').Value = '  +7.97%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.125'
$ws.Range('E10').Value = '  +13.30%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '41.86'
$ws.Range('E11').Value = '  +8.02%  '
$ws.Range('E12').Value = '  -0.72%  '
$ws.Range('D13').Value = '3.905.99'
$ws.Range('E13').Value = '  +1.47%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '8.46'
$ws.Range('E14').Value = '  +3.84%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '19.55'
$ws.Range('E15').Value = '  +2.99%  '
$ws.Range('D16').Value = '3.362.14'
$ws.Range('E16').Value = '  +1.57%  '
$ws.Range('B17').Value = 'Uniswap'
$ws.Range('C17').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '11.34'
$ws.Range('E17').Value = '  +7.93%  '
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').Value = '60.759.31'
$ws.Range('E18').Value = '  -0.50%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.02'
$ws.Range('E19').Value = '  +3.03%  '
$ws.Range('E20').Value = '  +19.15%  '
$ws.Range('E21').Value = '  +0.92%  '
$ws.Range('E22').Value = '  +12.81%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '13.14'
$ws.Range('E23').Value = '  +6.57%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '305.75'
$ws.Range('E24').Value = '  +3.66%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.12'
$ws.Range('E25').Value = '  +2.14%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '4.74'
$ws.Range('E26').Value = '  +5.56%  '
$ws.Range('E27').Value = '  +13.63%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '29.44'
$ws.Range('E28').Value = '  +2.66%  '
$ws.Range('E29').Value = '  +0.92%  '
$ws.Range('E30').Value = '  +1.28%  '
$ws.Range('E31').Value = '  +6.12%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '11.77'
$ws.Range('E32').Value = '  +5.91%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.60'
$ws.Range('E33').Value = '  +6.80%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '42.28'
$ws.Range('E34').Value = '  +6.38%  '
$ws.Range('E35').Value = '  +0.10%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0483'
$ws.Range('E36').Value = '  +1.35%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '52.27'
$ws.Range('E37').Value = '  -0.22%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.997'
$ws.Range('E38').Value = '  -0.18%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.41'
$ws.Range('E39').Value = '  +3.97%  '
$ws.Range('E40').Value = '  -2.52%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.02'
$ws.Range('E41').Value = '  +7.93%  '
$ws.Range('E42').Value = '  +4.85%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '136.37'
$ws.Range('E43').Value = '  +0.78%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.94'
$ws.Range('E44').Value = '  +4.94%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '16.87'
$ws.Range('E45').Value = '  +4.72%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.282'
$ws.Range('E46').Value = '  +0.44%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.22'
$ws.Range('E47').Value = '  +1.44%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '21.65'
$ws.Range('E48').Value = '  +4.02%  '
$ws.Range('D49').Value = '2.130.68'
$ws.Range('E49').Value = '  +0.90%  '
$ws.Range('D50').Value = '3.700.42'
$ws.Range('E50').Value = '  +1.40%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.36'
$ws.Range('E51').Value = '  +1.62%  '
